# "Add files via upload" — new form responses (Guilherme Luiz Maia Pinto,
# RA 501) were appended to the "Respostas ao formulário 1" sheet, and the
# e-mail-address cells (both the new ones and the pre-existing ones in
# column D) were turned into real mailto: hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New rows 8-11: four identical "Guilherme Luiz Maia Pinto" responses
# ---------------------------------------------------------------------
$guilhermeName  = "Guilherme Luiz Maia Pinto"
$guilhermeEmail = "guilhermemaia2011@live.com"

8..11 | ForEach-Object {
    $row = $_
    $ws.Range("B$row").Value = 501
    $ws.Range("C$row").Value = $guilhermeName
    $ws.Range("D$row").Value = $guilhermeEmail
    $ws.Hyperlinks.Add($ws.Range("D$row"), "mailto:$guilhermeEmail")
    $ws.Rows("$row").RowHeight = 15.75
}

# ---------------------------------------------------------------------
# Re-apply / add mailto hyperlinks to the pre-existing e-mail cells
# (D2, D5 = melaroccad@gmail.com ; D3, D6 = Bruno_clapis@hotmail.com,
# now stored without the trailing space it used to have)
# ---------------------------------------------------------------------
$maEmail = "melaroccad@gmail.com"
$ws.Range("D2").Value = $maEmail
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:$maEmail")

$brunoEmail = "Bruno_clapis@hotmail.com"
$ws.Range("D3").Value = $brunoEmail
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:$brunoEmail")

$ws.Range("D5").Value = $maEmail
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:$maEmail")

$ws.Range("D6").Value = $brunoEmail
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:$brunoEmail")

# ---------------------------------------------------------------------
# Stray formatted-but-empty cell further down the sheet (row 15), which
# is what stretches the sheet's dimension/used-range to A1:F15
# ---------------------------------------------------------------------
$ws.Range("D15").Font.Underline = 2

# Leave the selection where the user ended up after entering the data
$ws.Range("D7").Select()
